$wb = $excel.ActiveWorkbook

# Fix sheet name error: the 2nd and 3rd sheets were misnamed/shifted.
# Rename sheet 3 first (frees up the "SortEmployeeTableColumns" name),
# then rename sheet 2 to that freed-up name.
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws3.Name = "ShowEmployeesPerPage"
$ws2.Name = "SortEmployeeTableColumns"

# Update each sheet's remembered selection, and make the (renamed)
# "SortEmployeeTableColumns" sheet the active tab (it ends up selected
# last).
$ws3.Range("A25").Select()
$ws2.Range("B23").Select()
